$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44377
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 40000
$ws.Range("O2").Value = 40000
$ws.Range("P2").Value = 40000
$ws.Range("S2").Value = 2222

# Row 3
$ws.Range("D3").Value = 44658
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 28000
$ws.Range("O3").Value = 28000
$ws.Range("P3").Value = 28000
$ws.Range("S3").Value = 1556

# Row 4
$ws.Range("D4").Value = 44442
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 35000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = 35000
$ws.Range("R4").Value = 'Perú'
$ws.Range("S4").Value = 1944

# Row 5
$ws.Range("D5").Value = 44434
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("S5").Value = 1944

# Row 6
$ws.Range("D6").Value = 44435
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 35000
$ws.Range("O6").Value = 35000
$ws.Range("P6").Value = 35000
$ws.Range("R6").Value = 'Perú'
$ws.Range("S6").Value = 1944

# Row 7
$ws.Range("D7").Value = 44435
$ws.Range("M7").Value = 105
$ws.Range("R7").Value = 'Región de Arica y Parinacota'

# Row 8
$ws.Range("D8").Value = 44664
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("R8").Value = 'Perú'
$ws.Range("S8").Value = 1667

# Row 9
$ws.Range("D9").Value = 44662
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("S9").Value = 1667

# Row 10
$ws.Range("D10").Value = 44438
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 35000
$ws.Range("O10").Value = 35000
$ws.Range("P10").Value = 35000
$ws.Range("S10").Value = 1944

# Row 12
$ws.Range("D12").Value = 44629
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 35000
$ws.Range("O12").Value = 35000
$ws.Range("P12").Value = 35000
$ws.Range("S12").Value = 1944

# Row 13
$ws.Range("D13").Value = 44679
$ws.Range("M13").Value = 35
$ws.Range("R13").Value = 'Perú'

# Row 14
$ws.Range("D14").Value = 44679
$ws.Range("M14").Value = 55
$ws.Range("N14").Value = 28000
$ws.Range("O14").Value = 28000
$ws.Range("P14").Value = 28000
$ws.Range("S14").Value = 1556

# Row 15
$ws.Range("D15").Value = 44676
$ws.Range("M15").Value = 55
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 30000
$ws.Range("P15").Value = 28909
$ws.Range("S15").Value = 1606

# Row 16
$ws.Range("D16").Value = 44369
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 35000
$ws.Range("O16").Value = 35000
$ws.Range("P16").Value = 35000
$ws.Range("S16").Value = 1944

# Row 17
$ws.Range("D17").Value = 44418
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 35000
$ws.Range("O17").Value = 35000
$ws.Range("P17").Value = 35000
$ws.Range("S17").Value = 1944

# Row 18
$ws.Range("D18").Value = 44385
$ws.Range("M18").Value = 200
$ws.Range("N18").Value = 1200
$ws.Range("O18").Value = 1200
$ws.Range("P18").Value = 1200
$ws.Range("Q18").Value = '$/kilo'
$ws.Range("R18").Value = 'Perú'
$ws.Range("S18").Value = 1200
$ws.Range("T18").Value = 1

# Row 19
$ws.Range("D19").Value = 44405
$ws.Range("M19").Value = 10
$ws.Range("N19").Value = 35000
$ws.Range("O19").Value = 35000
$ws.Range("P19").Value = 35000
$ws.Range("Q19").Value = '$/caja 18 kilos'
$ws.Range("R19").Value = 'Región de Arica y Parinacota'
$ws.Range("S19").Value = 1944
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("D20").Value = 44690
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = 34000
$ws.Range("O20").Value = 34000
$ws.Range("P20").Value = 34000
$ws.Range("S20").Value = 1889

# Row 21
$ws.Range("D21").Value = 44357
$ws.Range("M21").Value = 10
$ws.Range("N21").Value = 38000
$ws.Range("O21").Value = 38000
$ws.Range("P21").Value = 38000
$ws.Range("R21").Value = 'Perú'
$ws.Range("S21").Value = 2111

# Row 22
$ws.Range("D22").Value = 44279
$ws.Range("M22").Value = 30
$ws.Range("O22").Value = 36000
$ws.Range("P22").Value = 35667
$ws.Range("S22").Value = 1982

# Row 23
$ws.Range("D23").Value = 44424
$ws.Range("M23").Value = 15
$ws.Range("O23").Value = 35000
$ws.Range("P23").Value = 35000
$ws.Range("S23").Value = 1944

# Row 24
$ws.Range("D24").Value = 44448
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 38000
$ws.Range("O24").Value = 38000
$ws.Range("P24").Value = 38000
$ws.Range("S24").Value = 2111

# Row 25
$ws.Range("D25").Value = 44645
$ws.Range("M25").Value = 5
$ws.Range("N25").Value = 30000
$ws.Range("O25").Value = 30000
$ws.Range("P25").Value = 30000
$ws.Range("S25").Value = 1667

# Row 26
$ws.Range("D26").Value = 44392
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 35000
$ws.Range("O26").Value = 35000
$ws.Range("P26").Value = 35000
$ws.Range("R26").Value = 'Región de Arica y Parinacota'
$ws.Range("S26").Value = 1944

# Row 27
$ws.Range("D27").Value = 44264
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 40000
$ws.Range("O27").Value = 40000
$ws.Range("P27").Value = 40000
$ws.Range("Q27").Value = '$/caja 18 kilos'
$ws.Range("S27").Value = 2222
$ws.Range("T27").Value = 18

# Row 28
$ws.Range("D28").Value = 44379
$ws.Range("N28").Value = 30000
$ws.Range("O28").Value = 30000
$ws.Range("P28").Value = 30000
$ws.Range("S28").Value = 1667

# Row 29
$ws.Range("D29").Value = 44294
$ws.Range("M29").Value = 15
$ws.Range("R29").Value = 'Región de Arica y Parinacota'

# Row 30
$ws.Range("D30").Value = 44431
$ws.Range("M30").Value = 30

# Row 31
$ws.Range("D31").Value = 44364
$ws.Range("M31").Value = 90
$ws.Range("N31").Value = 1700
$ws.Range("O31").Value = 1700
$ws.Range("P31").Value = 1700
$ws.Range("Q31").Value = '$/kilo'
$ws.Range("R31").Value = 'Región de Arica y Parinacota'
$ws.Range("S31").Value = 1700
$ws.Range("T31").Value = 1

# Row 32
$ws.Range("D32").Value = 44634
$ws.Range("N32").Value = 45000
$ws.Range("O32").Value = 45000
$ws.Range("P32").Value = 45000
$ws.Range("S32").Value = 2500

# Row 33
$ws.Range("D33").Value = 44669
$ws.Range("M33").Value = 40
$ws.Range("N33").Value = 32000
$ws.Range("O33").Value = 32000
$ws.Range("P33").Value = 32000
$ws.Range("S33").Value = 1778

# Row 34
$ws.Range("D34").Value = 44433
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = 35000
$ws.Range("O34").Value = 35000
$ws.Range("P34").Value = 35000
$ws.Range("S34").Value = 1944

# Row 35
$ws.Range("D35").Value = 44432
$ws.Range("M35").Value = 10
$ws.Range("R35").Value = 'Perú'

# Row 36
$ws.Range("D36").Value = 44449
$ws.Range("M36").Value = 20
$ws.Range("R36").Value = 'Región de Arica y Parinacota'

# Row 37
$ws.Range("D37").Value = 44671
$ws.Range("M37").Value = 20
$ws.Range("N37").Value = 32000
$ws.Range("O37").Value = 32000
$ws.Range("P37").Value = 32000
$ws.Range("R37").Value = 'Región de Arica y Parinacota'
$ws.Range("S37").Value = 1778
